# Apply the upload edit: update a handful of shared-string-backed cell
# values, widen columns F and G, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content edits -----------------------------------------------
# Row 5: "Easd" -> "Easd Rre"
$ws.Cells.Item(5, 6).Value = "Easd Rre"

# Row 4: "Tasd" -> "Tasd Mjh"
$ws.Cells.Item(4, 6).Value = "Tasd Mjh"

# Row 7: "Msad" -> "Msad Udfd"
$ws.Cells.Item(7, 7).Value = "Msad Udfd"

# Row 9: "Masda" -> "Masda Ouas"
$ws.Cells.Item(9, 7).Value = "Masda Ouas"

# --- Column widths for F and G -----------------------------------------
$ws.Columns.Item(6).ColumnWidth = 12.83
$ws.Columns.Item(7).ColumnWidth = 11.83

# --- Selection ----------------------------------------------------------
$ws.Range("G10").Select()
